# Generate Report for Handoff
# -----------------------------------------------------------------------
# The localization pipeline re-ran for a new source commit. The old
# commit's synthetic id "830954f8-9fd3-49ce-96f5-15d7d3de2793" is
# replaced everywhere by the new one "a2c8fd2a-d6da-4d47-9e2a-ac9da0ea9a20",
# timestamps advance, and because a fresh handoff was just generated the
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns for each language sheet reset to "not handed back
# yet" (blank / epoch date) and the obsolete "Latest Target File"
# hyperlink goes away.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$oldId = "830954f8-9fd3-49ce-96f5-15d7d3de2793"
$newId = "a2c8fd2a-d6da-4d47-9e2a-ac9da0ea9a20"

function Set-HyperlinkDisplay($ws, [string]$cellAddr, [string]$text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $cellAddr) {
            $hl.TextToDisplay = $text
        }
    }
}

function Remove-Hyperlink($ws, [string]$cellAddr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $cellAddr) {
            $hl.Delete()
        }
    }
}

# ------------------------------------------------------------------
# Overview sheet
# ------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newId.md"

$ov.Range("B2").Value = "e2e\$newId.md"
Set-HyperlinkDisplay $ov '$B$2' "e2e\$newId.md"

$ov.Range("G2").Value = "2016-08-18 21:01:20"

# ------------------------------------------------------------------
# zh-cn sheet
# ------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "$newId.md"
Set-HyperlinkDisplay $zh '$A$2' "$newId.md"

$zh.Range("G2").Value = "$newId.e3b1c021495b235c084b0a5b942d5acdda0fd047.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-18 21:01:14"

# A fresh handoff was just produced -- nothing has come back from
# localization yet, so target/handback file + handback datetime reset.
Remove-Hyperlink $zh '$I$2'
$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"

$zh.Range("J2").Value = ""

$zh.Range("K2").Value = "0001-01-01 00:00:00"

# ------------------------------------------------------------------
# de-de sheet
# ------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "$newId.md"
Set-HyperlinkDisplay $de '$A$2' "$newId.md"

$de.Range("G2").Value = "$newId.e3b1c021495b235c084b0a5b942d5acdda0fd047.de-de.xlf"
$de.Range("H2").Value = "2016-08-18 21:01:20"

Remove-Hyperlink $de '$I$2'
$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"

$de.Range("J2").Value = ""

$de.Range("K2").Value = "0001-01-01 00:00:00"

# ------------------------------------------------------------------
# Column widths on the language sheets shrink now that the target
# file / handback file columns are empty.
# ------------------------------------------------------------------
$zh.Columns.Item(9).AutoFit() | Out-Null
$zh.Columns.Item(10).AutoFit() | Out-Null
$de.Columns.Item(9).AutoFit() | Out-Null
$de.Columns.Item(10).AutoFit() | Out-Null
